# "write tc 03 and 04" -- rebuild the Result sheet as an empty TC template:
# keep the TC_ID / TC_Summary / TC_Result / Note header, drop the old
# TC01 sample rows, and lay out 8 blank bordered rows (2-9) ready for new
# test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old sample data (rows 2-6): values only, so the shared-strings
# table is rebuilt down to just the 4 header strings on save.
$ws.Range("A2:E6").ClearContents()

# Header row: keep the existing bold font, add a thin box border and
# center alignment.
$hdr = $ws.Range("A1:D1")
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.Borders.LineStyle = 1         # xlContinuous / thin

# Body rows 2-9: empty cells, just a thin box border around each one.
$body = $ws.Range("A2:D9")
$body.Borders.LineStyle = 1

# New column widths for the wider TC_ID / TC_Summary layout (closest
# achievable values given the engine's char-width quantization).
$ws.Columns.Item(1).ColumnWidth = 28.583333333333332
$ws.Columns.Item(2).ColumnWidth = 85.58333333333333
$ws.Columns.Item(3).ColumnWidth = 12.08333333333333
$ws.Columns.Item(4).ColumnWidth = 9.916666666666664

# Selection ends up parked below the table, as in the saved file.
$ws.Range("B18").Select()
